$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Map of cell address -> new literal text value (preserves exact formatting,
# e.g. trailing zeros / percent signs, matching the inlineStr cells in the sheet).
$updates = [ordered]@{
    "D2" = "304.50"
    "E2" = "4.06%"
    "D3" = "35.77"
    "E3" = "14.86%"
    "D4" = "5.095"
    "E4" = "2.77%"
    "D5" = "0.07843"
    "E5" = "4.86%"
    "D6" = "2.259"
    "E6" = "-1.05%"
    "D7" = "8.119"
    "E7" = "4.23%"
    "D8" = "4.005"
    "E8" = "6.07%"
    "D9" = "0.9268"
    "E9" = "0.88%"
    "D10" = "0.09854"
    "E10" = "5.80%"
    "D11" = "0.1824"
    "E11" = "5.42%"
    "D12" = "0.08764"
    "E12" = "5.15%"
    "D13" = "0.03415"
    "E13" = "4.13%"
    "D14" = "0.09951"
    "E14" = "0.16%"
    "D15" = "0.001478"
    "E15" = "-1.77%"
    "D16" = "0.005779"
    "E16" = "1.37%"
    "D17" = "3.483"
    "E17" = "0.18%"
    "E18" = "-1.23%"
    "D19" = "0.3432"
    "E19" = "2.94%"
    "D20" = "0.1320"
    "E20" = "0.45%"
    "D21" = "4.549"
    "E21" = "10.98%"
    "E22" = "6.48%"
    "E23" = "3.14%"
    "D24" = "0.001241"
    "E24" = "1.78%"
    "D25" = "0.004504"
    "E25" = "4.61%"
    "D26" = "0.0001301"
    "E26" = "0.16%"
    "D27" = "0.0002701"
    "E27" = "-20.37%"
    "D39" = "0.01755"
    "E39" = "8.22%"
    "E40" = "3.01%"
    "D41" = "0.008016"
    "E41" = "8.01%"
    "D42" = "0.1424"
    "E42" = "4.72%"
    "D43" = "0.008488"
    "E43" = "-13.76%"
    "D44" = "0.002213"
    "E44" = "-0.26%"
    "D45" = "0.009132"
    "D46" = "0.00006163"
    "E46" = "1.04%"
    "D47" = "0.00000000751"
    "E47" = "0.15%"
    "D48" = "4.049"
    "E48" = "58.73%"
    "D49" = "0.002692"
    "E49" = "34.70%"
    "D50" = "0.00002102"
    "E50" = "0.15%"
    "D51" = "0.0002002"
    "E51" = "0.15%"
}

foreach ($addr in $updates.Keys) {
    $cell = $ws.Range($addr)
    # Force text number format so Excel keeps the string verbatim instead of
    # re-interpreting it as a number/percentage (which would drop trailing zeros
    # or rescale percent values).
    $cell.NumberFormat = "@"
    $cell.Value = $updates[$addr]
}
